$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "96.786.97"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.43%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.644.04"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("E4").Value = "  +0.00%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "242.14"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.37%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.87"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +19.93%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "655.61"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.19%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.423"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +4.40%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "1.07"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +3.10%  "
$ws.Range("E10").Value = "  -0.05%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "3.640.53"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +1.48%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "44.22"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +2.41%  "
$ws.Range("E13").Value = "  +1.18%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.48"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.33%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "4.331.49"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +1.75%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "96.509.19"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("E17").Value = "  +0.45%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "3.626.40"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +1.68%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "7.86"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.26%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "12.84"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.98%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "18.29"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +2.81%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.539"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +9.20%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "512.26"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("B24").Value = "SuiNetwork"
$ws.Range("C24").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "3.44"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.36%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.0000206"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +2.87%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "6.89"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.64%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "101.24"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +4.88%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "13.00"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +1.59%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.171"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +15.25%  "
$ws.Range("E30").Value = "  +1.49%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "11.82"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +3.42%  "
$ws.Range("E32").Value = "  +0.03%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.186"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +1.21%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "32.89"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +3.93%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("E36").Value = "  +7.40%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.582"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +2.74%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "8.80"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.88%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "615.78"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.43%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "41.44"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +20.70%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.159"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +5.52%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.949"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +4.33%  "
$ws.Range("E43").Value = "  +6.18%  "
$ws.Range("E44").Value = "  -0.01%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "6.14"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +7.84%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0444"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +6.47%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.422"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +23.32%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.31"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.41%  "
$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "23.61"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("E50").Value = "  +5.15%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "54.36"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +2.11%  "
